$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/20/2025  Through  10/26/2025"

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("N14").Value = -88.888888888888
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 22
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 175
$ws.Range("M15").Value = 214.285714285714
$ws.Range("N15").Value = 15.789473684210
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 0
$ws.Range("L16").Value = -23.622047244094
$ws.Range("M16").Value = -27.611940298507
$ws.Range("N16").Value = -80.244399185336
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -83.333333333333
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -52.631578947368
$ws.Range("I17").Value = 184
$ws.Range("J17").Value = 204
$ws.Range("K17").Value = -9.803921568627
$ws.Range("L17").Value = 8.235294117647
$ws.Range("M17").Value = 85.858585858585
$ws.Range("N17").Value = -19.650655021834
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 111
$ws.Range("J18").Value = 138
$ws.Range("K18").Value = -19.565217391304
$ws.Range("L18").Value = -5.932203389830
$ws.Range("M18").Value = -18.382352941176
$ws.Range("N18").Value = -88.425443169968
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -9.090909090909
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = -30.909090909090
$ws.Range("I19").Value = 468
$ws.Range("J19").Value = 630
$ws.Range("K19").Value = -25.714285714285
$ws.Range("L19").Value = -16.873889875666
$ws.Range("M19").Value = 67.741935483871
$ws.Range("N19").Value = 23.482849604221
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -30
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = -43.333333333333
$ws.Range("I20").Value = 226
$ws.Range("J20").Value = 280
$ws.Range("K20").Value = -19.285714285714
$ws.Range("L20").Value = -20.979020979021
$ws.Range("M20").Value = 119.417475728155
$ws.Range("N20").Value = -85.033112582781
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -39.393939393939
$ws.Range("F21").Value = 79
$ws.Range("G21").Value = 126
$ws.Range("H21").Value = -37.301587301587
$ws.Range("I21").Value = 1109
$ws.Range("J21").Value = 1380
$ws.Range("K21").Value = -19.637681159420
$ws.Range("L21").Value = -13.155833985904
$ws.Range("M21").Value = 45.347313237221
$ws.Range("N21").Value = -69.160177975528
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("J23").Value = 55
$ws.Range("K23").Value = -18.181818181818
$ws.Range("L23").Value = -34.782608695652
$ws.Range("M23").Value = 25
$ws.Range("C24").Value = 41
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 64
$ws.Range("F24").Value = 106
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = 8.163265306122
$ws.Range("I24").Value = 963
$ws.Range("J24").Value = 915
$ws.Range("K24").Value = 5.245901639344
$ws.Range("L24").Value = 2.884615384615
$ws.Range("M24").Value = 43.090638930163
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -28.571428571428
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = -8
$ws.Range("I25").Value = 269
$ws.Range("J25").Value = 344
$ws.Range("K25").Value = -21.802325581395
$ws.Range("L25").Value = -27.882037533512
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -7.692307692307
$ws.Range("F26").Value = 28
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = -26.315789473684
$ws.Range("I26").Value = 317
$ws.Range("J26").Value = 295
$ws.Range("K26").Value = 7.457627118644
$ws.Range("L26").Value = 23.346303501945
$ws.Range("M26").Value = -2.760736196319
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 26
$ws.Range("K27").Value = 36.842105263157
$ws.Range("L27").Value = 18.181818181818
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -80
$ws.Range("I28").Value = 34
$ws.Range("J28").Value = 42
$ws.Range("K28").Value = -19.047619047619
$ws.Range("L28").Value = 17.241379310344

# --- Cells converting from Number to Text (style -> 13, shared text) ---
# Force text entry via NumberFormat "@", then restore the exact donor style
# via PasteSpecial(formats) so the saved style index matches the original
# document's deduplicated style table precisely.
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null

# --- Cells converting from Text to Number (style -> 14/15, numeric) ---
$ws.Range("C15").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("F15").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("F15").PasteSpecial(-4122) | Out-Null
$ws.Range("D18").Value = 6
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = -83.333333333333
$ws.Range("L14").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("C27").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("F27").PasteSpecial(-4122) | Out-Null
$ws.Range("C28").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
